$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Title paragraph ("กิตติกรรมประกาศ"): bump font size 22pt -> 24pt
#    (w:sz/w:szCs 44 -> 48) on both the paragraph mark run props and the
#    text run itself.
# ---------------------------------------------------------------------
$title = $d.Paragraphs(1).Range
$title.Font.Size = 24
$title.Font.SizeBi = 24

# ---------------------------------------------------------------------
# 2) Fourth body paragraph: insert the missing "น" between "ท่า" and
#    "ในทีม " (ท่า -> ท่าน), which splits the original run into three
#    runs once Word (re)applies direct character formatting to the
#    inserted text and its immediate neighbour.
# ---------------------------------------------------------------------
$full = $d.Content.Text
$splitPos = $full.IndexOf("ในทีม")

# Insert the new character right before "ในทีม"
$insertion = $d.Range($splitPos, $splitPos)
$insertion.InsertAfter([char]0x0E19)

# The inserted "น" now occupies [$splitPos, $splitPos+1); "ในทีม " (with
# the trailing space) occupies [$splitPos+1, $splitPos+7). Forcing a
# (no-op) font-size change on each sub-range splits the run boundaries
# to match, then restores the original 16pt size.
$nRange = $d.Range($splitPos, $splitPos + 1)
$nRange.Font.Size = 99
$nRange2 = $d.Range($splitPos, $splitPos + 1)
$nRange2.Font.Size = 16

$tailRange = $d.Range($splitPos + 1, $splitPos + 7)
$tailRange.Font.Size = 99
$tailRange2 = $d.Range($splitPos + 1, $splitPos + 7)
$tailRange2.Font.Size = 16

# ---------------------------------------------------------------------
# 3) Signature line: drop the "นางสาว" honorific prefix.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("นางสาวพิชาดา เลิศประเสริฐกิจ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "พิชาดา เลิศประเสริฐกิจ", 2)
